# Task planning update 1.2
# Update the "Status" column (F) of Sheet1 for tasks whose progress moved
# forward, then move the sheet's selection to where work is now focused
# (row 25, which is scrolled into view as the new top area).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F15").Value = "Done"
$ws.Range("F16").Value = "In progress"
$ws.Range("F17").Value = "In progress"
$ws.Range("F18").Value = "In progress"
$ws.Range("F19").Value = "Done"
$ws.Range("F20").Value = "Done"
$ws.Range("F21").Value = "In progress"
$ws.Range("F22").Value = "In progress"
$ws.Range("F23").Value = "Done"
$ws.Range("F24").Value = "Done"
$ws.Range("F25").Value = "Done"
$ws.Range("F26").Value = "Done"
$ws.Range("F27").Value = "Done"
$ws.Range("F28").Value = "Done"
$ws.Range("F29").Value = "Done"
$ws.Range("F30").Value = "Done"
$ws.Range("F31").Value = "Done"
$ws.Range("F32").Value = "Done"
$ws.Range("F33").Value = "In progress"
$ws.Range("F34").Value = "Done"
$ws.Range("F35").Value = "Done"
$ws.Range("F36").Value = "In progress"
$ws.Range("F37").Value = "Done"
$ws.Range("F38").Value = "In progress"
$ws.Range("F41").Value = "Done"
$ws.Range("F42").Value = "Done"
$ws.Range("F51").Value = "Done"
$ws.Range("F53").Value = "Done"

# Move the active selection to reflect where the user is now working.
$null = $ws.Range("F25").Select()
